# Orginfo, Excel einlesen, Flask bestandteil von Main
#
# Adds 9 rows of "s" placeholder values into column A (rows 3-11) of the
# "Daten" worksheet, used as simple read markers for the Excel-import logic.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Daten")

for ($row = 3; $row -le 11; $row++) {
    $ws.Cells.Item($row, 1).Value = "s"
}

$ws.Range("A12").Select()
